$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (Miscellaneous / Adrenal adenoma / Clip 1 B-mode): add new YouTube link in D13
$ws.Range("D13").Value = "https://youtu.be/xBfd04F4Ni8 "

# Row 8 (Liver / HNF1a-mutated hepatocellular adenoma - Hyperechoic / Clip 1 B-mode + Color):
# the previously-unlinked YouTube text in D8 becomes an actual hyperlink (same URL, trailing space added)
$ws.Range("D8").Value = "https://youtu.be/91M82AIMyu0 "

# Turn both new cells into real hyperlinks
$ws.Hyperlinks.Add($ws.Range("D13"), "https://youtu.be/xBfd04F4Ni8")
$ws.Hyperlinks.Add($ws.Range("D8"), "https://youtu.be/91M82AIMyu0")

# Make sure both cells use the same "Hyperlink" cell style as the other link cells (D3, D9)
$ws.Range("D13").Style = "Collegamento ipertestuale"
$ws.Range("D8").Style = "Collegamento ipertestuale"

# Update the saved cursor/selection position
$ws.Range("D14").Select()
